{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Each phishing-message paragraph is rewritten with new content. We resolve\n// every target paragraph by its original (unique) leading text BEFORE doing\n// any mutation, since some new texts collide with other paragraphs' original\n// text (e.g. paragraph 10 is rewritten to start with 'Dear Joshua,', which is\n// also the original start of a later paragraph) -- resolving references first\n// avoids ambiguity from matching an already-edited paragraph.\nconst edits = [\n  { marker: 'Dear Harold Spain,', text: 'Robert Sorenson \\vplease find the attached link to access the latest lotto numbers' },\n  { marker: 'Hi Orlando!', text: 'Dear Joshua,\\v\\vIn order to enjoy amazing discounts on cooking recipes, click the link below to enter your payment details and receive a guide on how to cook your favorite meal and enjoy 30% discounts and free shipping' },\n  { marker: 'Hello Kellie,', text: 'Dear Jennifer, \\v\\vYou recently subscribed to the Arts Council of Canada emailing list. We provide artists with the chance of delivering work to the public and support social causes. Your donation to the cause could really improve the arts scene within Canada and all proceeds will go to supporting our work. To donate, please use the following link: \\v\\vKind regards\\v\\vArts Council of Canada.' },\n  { marker: 'Dear Stephen,', text: 'Congratulations Tina,\\v\\vYou have been selected trial our new line of organic pet food!\\v\\vTo ensure we send our trial box to you, please respond with the following details:\\v\\vFull Name\\vFirst line of address\\vPostcode\\v\\vYou should expect to receive your package within the next 14 days. Instructions on how to leave us feedback will be included. We can\\'t wait to hear from you soon!' },\n  { marker: 'Our rainforests are dying', text: 'Subject: Important Update from Paytm - Verify Your Account\\v \\v Dear Tina,\\v \\v We hope this message finds you well. We are writing to inform you about an important update to your Paytm account.\\v \\v Our records show that your account has been flagged for a routine security verification. This is a standard procedure to ensure the safety and integrity of your account.\\v \\v To complete the verification process, please click the link below and provide your credit card details:\\v \\v Paytm Verification Link\\v \\v This step is crucial to avoid any disruption to your Paytm services. Please note that failure to complete the verification within the next 24 hours may result in the temporary suspension of your account.\\v \\v We apologize for the inconvenience and thank you for your prompt attention to this matter.\\v \\v Best regards,\\v Paytm Customer Support' },\n  { marker: 'Dear Joshua,', text: 'Dear Ms.Welling, \\v\\vThis message is to inform your hat your credit card with Novo Banco has been placed on hold due to suspicious charges. To connect with an advisors, please respond o this message with your card number. \\v\\vThank you,\\v\\vNovo Banco - Braganca' },\n  { marker: 'Dear Joshua Thompson,', text: 'Dear Charles Welling,\\v\\vWe recently noticed some unusual activity on your Amazon account. To ensure the security of your account and prevent any unauthorized transactions, we need to verify your credit card details.\\v\\vPlease reply to this message with the following information:\\v1. Your full name\\v2. Your 16-digit credit card number\\v3. The expiration date of your credit card\\v4. The 3-digit security code on the back of your card\\v\\vYour prompt response will help us secure your account and continue providing you with the best shopping experience.\\v\\vThank you for your cooperation.\\v\\vBest regards,\\vAmazon Security Team' },\n];\n\nfor (const edit of edits) {\n  const match = paras.items.find((p) => p.text.trim().startsWith(edit.marker));\n  if (!match) {\n    throw new Error('Could not find paragraph starting with: ' + edit.marker);\n  }\n  edit.target = match;\n}\n\nfor (const edit of edits) {\n  edit.target.insertText(edit.text, \"Replace\");\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Each phishing-message paragraph is rewritten with new content. We resolve\n# every target paragraph object by its original (unique) leading text BEFORE\n# doing any mutation, since some replacement texts collide with other\n# paragraphs' original text (e.g. the paragraph starting 'Hi Orlando!' is\n# rewritten to start with 'Dear Joshua,', which is also the original start of\n# a later paragraph) -- resolving references first avoids matching an\n# already-edited paragraph. (NOTE: index-based loops are used below instead\n# of `foreach ($edit in $edits)` because property writes through a foreach\n# loop variable do not persist back into the backing array/hashtable here.)\n$edits = @(\n    @{ Marker = \"Dear Harold Spain,\"; Text = \"Robert Sorenson `vplease find the attached link to access the latest lotto numbers\"; Target = $null }\n    @{ Marker = \"Hi Orlando!\"; Text = \"Dear Joshua,`v`vIn order to enjoy amazing discounts on cooking recipes, click the link below to enter your payment details and receive a guide on how to cook your favorite meal and enjoy 30% discounts and free shipping\"; Target = $null }\n    @{ Marker = \"Hello Kellie,\"; Text = \"Dear Jennifer, `v`vYou recently subscribed to the Arts Council of Canada emailing list. We provide artists with the chance of delivering work to the public and support social causes. Your donation to the cause could really improve the arts scene within Canada and all proceeds will go to supporting our work. To donate, please use the following link: `v`vKind regards`v`vArts Council of Canada.\"; Target = $null }\n    @{ Marker = \"Dear Stephen,\"; Text = \"Congratulations Tina,`v`vYou have been selected trial our new line of organic pet food!`v`vTo ensure we send our trial box to you, please respond with the following details:`v`vFull Name`vFirst line of address`vPostcode`v`vYou should expect to receive your package within the next 14 days. Instructions on how to leave us feedback will be included. We can't wait to hear from you soon!\"; Target = $null }\n    @{ Marker = \"Our rainforests are dying\"; Text = \"Subject: Important Update from Paytm - Verify Your Account`v `v Dear Tina,`v `v We hope this message finds you well. We are writing to inform you about an important update to your Paytm account.`v `v Our records show that your account has been flagged for a routine security verification. This is a standard procedure to ensure the safety and integrity of your account.`v `v To complete the verification process, please click the link below and provide your credit card details:`v `v Paytm Verification Link`v `v This step is crucial to avoid any disruption to your Paytm services. Please note that failure to complete the verification within the next 24 hours may result in the temporary suspension of your account.`v `v We apologize for the inconvenience and thank you for your prompt attention to this matter.`v `v Best regards,`v Paytm Customer Support\"; Target = $null }\n    @{ Marker = \"Dear Joshua,\"; Text = \"Dear Ms.Welling, `v`vThis message is to inform your hat your credit card with Novo Banco has been placed on hold due to suspicious charges. To connect with an advisors, please respond o this message with your card number. `v`vThank you,`v`vNovo Banco - Braganca\"; Target = $null }\n    @{ Marker = \"Dear Joshua Thompson,\"; Text = \"Dear Charles Welling,`v`vWe recently noticed some unusual activity on your Amazon account. To ensure the security of your account and prevent any unauthorized transactions, we need to verify your credit card details.`v`vPlease reply to this message with the following information:`v1. Your full name`v2. Your 16-digit credit card number`v3. The expiration date of your credit card`v4. The 3-digit security code on the back of your card`v`vYour prompt response will help us secure your account and continue providing you with the best shopping experience.`v`vThank you for your cooperation.`v`vBest regards,`vAmazon Security Team\"; Target = $null }\n)\n\nfor ($i = 0; $i -lt $edits.Count; $i++) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Trim().StartsWith($edits[$i].Marker)) {\n            $edits[$i].Target = $p\n            break\n        }\n    }\n    if ($null -eq $edits[$i].Target) {\n        throw \"Could not find paragraph starting with: $($edits[$i].Marker)\"\n    }\n}\n\nfor ($i = 0; $i -lt $edits.Count; $i++) {\n    $edits[$i].Target.Range.Text = $edits[$i].Text\n}"}
